# series.xlsx: repair lag_ts/seq.period output table.
#   - drop the "code" column (was column B, all placeholder "x"/" " values)
#   - insert a blank separator row above the second group of series (row 6)
#   - leave the remaining data columns (former C:F) as the new B:E
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column B's formatting before removing it so the deletion doesn't
# leave a dangling style reference behind on the shifted column records.
$ws.Columns("B").ClearFormats()

# Insert a blank row before the "BD_WN___" group (old row 6), pushing the
# last three series down by one row.
$ws.Rows("6").Insert()

# Drop the old "code" column entirely; data that used to live in C:F now
# becomes B:E.
$ws.Columns("B").Delete()

# Match the author's last on-screen selection (whole row 6 highlighted).
$ws.Rows("6").EntireRow.Select() | Out-Null

Write-Host "Removed code column, inserted separator row, reselected row 6."
